$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the cell-formatting (style) of A15 (which carries the bordered/bold/centered style)
# onto the new A16 cell before setting its value.
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"

$ws.Range("C16").Value = 0.9801883432228856
$ws.Range("D16").Value = 1.059928533930848
$ws.Range("E16").Value = 0.9808058266454883
$ws.Range("F16").Value = 0.9801883432228856
$ws.Range("G16").Value = 1.032208412466209
$ws.Range("H16").Value = 0.9500791628242369
$ws.Range("I16").Value = 0.9790772385206749
$ws.Range("J16").Value = 1.059928533930848
$ws.Range("K16").Value = 1.020367180288168
$ws.Range("L16").Value = 1.000277761755527
$ws.Range("M16").Value = 0.9970479196017238
